# Auto-generated edit script applying the cryptos.xlsx diff via Excel COM interop.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cell updates (coin names, links, percent-change strings).
$ws.Range("D2").Value = "42.824.37"
$ws.Range("E2").Value = "  -5.05%  "
$ws.Range("D3").Value = "2.211.46"
$ws.Range("E3").Value = "  -6.26%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").Value = "  -7.34%  "
$ws.Range("E7").Value = "  -6.40%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -7.75%  "
$ws.Range("E10").Value = "  -9.51%  "
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("E12").Value = "  -9.60%  "
$ws.Range("E13").Value = "  -8.92%  "
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E15").Value = "  -11.57%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.547.17"
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("E17").Value = "  -6.65%  "
$ws.Range("D18").Value = "2.211.94"
$ws.Range("E18").Value = "  -5.90%  "
$ws.Range("D19").Value = "42.803.54"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("E21").Value = "  -9.29%  "
$ws.Range("E22").Value = "  -10.41%  "
$ws.Range("E23").Value = "  -10.66%  "
$ws.Range("E24").Value = "  -10.57%  "
$ws.Range("E25").Value = "  -8.82%  "
$ws.Range("E26").Value = "  -7.48%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E29").Value = "  -9.28%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("E31").Value = "  -12.38%  "
$ws.Range("E32").Value = "  -8.80%  "
$ws.Range("E33").Value = "  -7.74%  "
$ws.Range("E34").Value = "  -7.61%  "
$ws.Range("E35").Value = "  -8.08%  "
$ws.Range("E36").Value = "  -7.04%  "
$ws.Range("E37").Value = "  +6.66%  "
$ws.Range("E38").Value = "  -6.62%  "
$ws.Range("E39").Value = "  +8.04%  "
$ws.Range("E40").Value = "  -8.31%  "
$ws.Range("E41").Value = "  -5.28%  "
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("E43").Value = "  -7.73%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "1.812.70"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("E46").Value = "  -4.63%  "
$ws.Range("E47").Value = "  -12.17%  "
$ws.Range("E48").Value = "  -9.34%  "
$ws.Range("E49").Value = "  -5.56%  "
$ws.Range("E50").Value = "  -12.38%  "
$ws.Range("E51").Value = "  -9.52%  "

# Price cells whose new text looks like a plain number (e.g. "9.97").
# Excel auto-converts such literals to numeric values on assignment, so
# force each target cell to Text format first, assign the literal string,
# then restore the default "Normal" style so no visible formatting changes.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.589"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.560"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.859"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0889"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0325"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.54"
$ws.Range("D51").Style = "Normal"
